$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2009年" data row (row 2) is being removed entirely; the following
# rows ("2010年", "2011年") shift up by one row, and the sheet dimension
# shrinks accordingly (Excel/iron_native handles the dimension + shifting
# automatically as part of a real row deletion).
$ws.Rows(2).Delete()
